$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - update F column (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 4536
$ws1.Range("F13").Value = 586
$ws1.Range("F15").Value = 507
$ws1.Range("F18").Value = 1755
$ws1.Range("F19").Value = 1294
$ws1.Range("F21").Value = 1535
$ws1.Range("F31").Value = 3449
$ws1.Range("F34").Value = 226
$ws1.Range("F36").Value = 1672

# Sheet "全部类型" (all types) - update F column (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 4536
$ws4.Range("F13").Value = 586
$ws4.Range("F15").Value = 507
$ws4.Range("F19").Value = 1755
$ws4.Range("F20").Value = 1294
$ws4.Range("F22").Value = 1535
$ws4.Range("F32").Value = 3449
$ws4.Range("F36").Value = 226
$ws4.Range("F38").Value = 1672
